$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column O: "Mensagem" header + "User name already exists" value,
# mirroring the existing header style (bold) used by A1:N1.
$ws.Range("O1").Value = "Mensagem"
$ws.Range("O1").Font.Bold = $true
$ws.Range("O2").Value = "User name already exists"

# Best-effort approximation of the bestFit column width Excel would have
# computed for the new column (exact bestFit pixel width isn't reachable
# through this runtime's ColumnWidth setter, so we pick the closest we can).
$ws.Columns.Item(15).ColumnWidth = 22.6

# Restore the view/selection state recorded in the saved workbook.
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("L24:L25").Select() | Out-Null
